$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the "Adjusted by" (column H) for the comments that were
#     actioned by Mazrouaa ---
$ws.Range("H42").Value = "Mazrouaa"
$ws.Range("H45").Value = "Mazrouaa"
$ws.Range("H46").Value = "Mazrouaa"
$ws.Range("H48").Value = "Mazrouaa"
$ws.Range("H49").Value = "Mazrouaa"
$ws.Range("H50").Value = "Mazrouaa"
$ws.Range("H51").Value = "Mazrouaa"
$ws.Range("H52").Value = "Mazrouaa"

# --- Update the status (column I) for comment row 44 ---
$ws.Range("I44").Value = "In Progress"

# --- Update the saved scroll position / selection on the sheet ---
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
$ws.Range("I43").Select() | Out-Null
